$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Candidate / Status / Action Date values between row 6 and row 8
$row6D = $ws.Range("D6").Value2
$row6E = $ws.Range("E6").Value2
$row6F = $ws.Range("F6").Value2

$row8D = $ws.Range("D8").Value2
$row8E = $ws.Range("E8").Value2
$row8F = $ws.Range("F8").Value2

$ws.Range("D6").Value = $row8D
$ws.Range("E6").Value = $row8E
$ws.Range("F6").Value = $row8F

$ws.Range("D8").Value = $row6D
$ws.Range("E8").Value = $row6E
$ws.Range("F8").Value = $row6F
